$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.03377890586853
$ws.Range("B1").Value = 4.708751678466797
$ws.Range("C1").Value = 3.541271448135376
$ws.Range("D1").Value = 0.9001575112342834
$ws.Range("E1").Value = 0.4726637303829193
